$d = $word.ActiveDocument

# 1) Conta Bancária: add missing "-7" suffix to the account number
$d.Content.Find.Execute("Conta Bancária: 10.738", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Conta Bancária: 10.738-7", 2)

# 2) Vigência das Atividades do Projeto: append the activity validity period
$d.Content.Find.Execute("Vigência das Atividades do Projeto:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Vigência das Atividades do Projeto: 23/02/2021 a 23/02/2022", 2)
